$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.02453427314758301
$ws.Range("C2").Value = 0.04562606811523438
$ws.Range("D2").Value = 0.01350531578063965
$ws.Range("E2").Value = 0.03302760124206543
$ws.Range("F2").Value = 0.008483648300170898
$ws.Range("G2").Value = 0.09333915710449218
$ws.Range("H2").Value = 0.02755804061889649
$ws.Range("I2").Value = 0.03097348213195801
$ws.Range("J2").Value = 0.0206578254699707
$ws.Range("K2").Value = 0.02799925804138183
$ws.Range("L2").Value = 0.003416252136230469
$ws.Range("M2").Value = 0.02531528472900391
$ws.Range("B3").Value = 0.1226459980010986
$ws.Range("C3").Value = 0.04267120361328125
$ws.Range("D3").Value = 0.01536340713500977
$ws.Range("E3").Value = 0.01508159637451172
$ws.Range("F3").Value = 0.0113978385925293
$ws.Range("G3").Value = 0.007460880279541016
$ws.Range("H3").Value = 0.1285584926605225
$ws.Range("I3").Value = 0.04587903022766113
$ws.Range("J3").Value = 0.08203740119934082
$ws.Range("K3").Value = 0.02377519607543945
$ws.Range("L3").Value = 0.02572412490844727
$ws.Range("M3").Value = 0.01503868103027344
$ws.Range("B4").Value = 0.04863910675048828
$ws.Range("C4").Value = 0.02747330665588379
$ws.Range("D4").Value = 0.0225034236907959
$ws.Range("E4").Value = 0.01230783462524414
$ws.Range("F4").Value = 0.1012078762054443
$ws.Range("G4").Value = 0.01251845359802246
$ws.Range("H4").Value = 0.02721915245056152
$ws.Range("I4").Value = 0.0156346321105957
$ws.Range("J4").Value = 0.02772345542907715
$ws.Range("K4").Value = 0.02142224311828613
$ws.Range("L4").Value = 0.04239592552185058
$ws.Range("M4").Value = 0.009864234924316406
$ws.Range("B5").Value = 0.03178791999816895
$ws.Range("C5").Value = 0.02495737075805664
$ws.Range("D5").Value = 0.03640303611755371
$ws.Range("E5").Value = 0.03060135841369629
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.02167863845825195
$ws.Range("I5").Value = 0.03239674568176269
$ws.Range("J5").Value = 0.02532310485839844
$ws.Range("K5").Value = 0.03071327209472656
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("B6").Value = 0.6837720394134521
$ws.Range("C6").Value = 0.03015332221984863
$ws.Range("D6").Value = 0.6274021625518799
$ws.Range("E6").Value = 0.04070558547973633
$ws.Range("F6").Value = 1.200747680664062
$ws.Range("G6").Value = 0.02355718612670898
$ws.Range("H6").Value = 0.3341116428375244
$ws.Range("I6").Value = 0.02347135543823242
$ws.Range("J6").Value = 0.5930277824401855
$ws.Range("K6").Value = 0.02374272346496582
$ws.Range("L6").Value = 0.4998091697692871
$ws.Range("M6").Value = 0.01894192695617676
$ws.Range("B7").Value = 0.8560727596282959
$ws.Range("C7").Value = 0.1082107067108154
$ws.Range("D7").Value = 0.4238080024719239
$ws.Range("E7").Value = 0.05590958595275879
$ws.Range("F7").Value = 0.5734320640563965
$ws.Range("G7").Value = 0.02772893905639649
$ws.Range("H7").Value = 0.9348299980163575
$ws.Range("I7").Value = 0.1075291156768799
$ws.Range("J7").Value = 0.3289021015167236
$ws.Range("K7").Value = 0.04420270919799805
$ws.Range("L7").Value = 0.6750794887542725
$ws.Range("M7").Value = 0.02424759864807129
